# Update "南宁-漫展信息" workbook (gh-pages data refresh at commit 456a3b4)
# Sheet "展览" (index 1): the oldest entry (南宁·AP动漫游戏嘉年华) rolled off the
#   top of the list and a new entry (南宁·蔚蓝档案only) rolled off the bottom,
#   every remaining row's 想去人数/最低票价 (and other scraped fields) refreshed.
# Sheet "全部类型" (index 4): same kind of refresh/roll, with one additional
#   concert entry (久石让作品视听音乐会) inserted near the top.
# Sheets "演出" (index 2) and "本地生活" (index 3) are untouched by this commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (1st sheet): 14 data+header rows -> 12 data+header rows.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Drop the oldest record (old row 2) -- everything below shifts up one row.
$ws1.Rows.Item(2).Delete()
# Drop what is now the trailing record (old row 14, now row 13).
$ws1.Rows.Item(13).Delete()

# Renumber the "0"-based index column (A) so it stays sequential (1..11).
for ($r = 2; $r -le 12; $r++) {
    $ws1.Cells.Item($r, 1).Value = $r - 1
}

# Refresh every remaining row's scraped fields (B:I) to the latest values.
# (Column B keeps its YYYY-MM-DD text forced via NumberFormat "@" so Excel
# doesn't silently reinterpret it as a date serial; the style is restored
# to "Normal" right after so no stray number format is left on the cell.)
$ws1.Cells.Item(2,2).NumberFormat = "@"
$ws1.Cells.Item(2,2).Value = '2024-06-15'
$ws1.Cells.Item(2,2).Style = "Normal"
$ws1.Cells.Item(2,3).Value = '南宁·星STAR国潮嘉年华'
$ws1.Cells.Item(2,4).Value = '亭洪路45号 百益上河城'
$ws1.Cells.Item(2,5).Value = '2024.06.15 09:00-06.16 17:00'
$ws1.Cells.Item(2,6).Value = 73
$ws1.Cells.Item(2,7).Value = 50
$ws1.Cells.Item(2,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86198'
$ws1.Cells.Item(2,9).Value = '//i0.hdslb.com/bfs/openplatform/202405/orwMgait1716448294056.jpeg'
$ws1.Cells.Item(3,2).NumberFormat = "@"
$ws1.Cells.Item(3,2).Value = '2024-06-22'
$ws1.Cells.Item(3,2).Style = "Normal"
$ws1.Cells.Item(3,3).Value = '南宁·排球少年ONLY（取消）'
$ws1.Cells.Item(3,4).Value = '亭洪路45号 水明漾宴会中心'
$ws1.Cells.Item(3,5).Value = '2024.06.22 09:45-06.22 17:00'
$ws1.Cells.Item(3,6).Value = 63
$ws1.Cells.Item(3,7).Value = '不可售'
$ws1.Cells.Item(3,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86465'
$ws1.Cells.Item(3,9).Value = '//i0.hdslb.com/bfs/openplatform/202405/GaaD97dL1716883956953.jpeg'
$ws1.Cells.Item(4,2).NumberFormat = "@"
$ws1.Cells.Item(4,2).Value = '2024-07-06'
$ws1.Cells.Item(4,2).Style = "Normal"
$ws1.Cells.Item(4,3).Value = '南宁·小蜜蜂动漫嘉年华2.0'
$ws1.Cells.Item(4,4).Value = '亭洪路45号 百益上河城'
$ws1.Cells.Item(4,5).Value = '2024.07.06 10:00-07.06 17:00'
$ws1.Cells.Item(4,6).Value = 245
$ws1.Cells.Item(4,7).Value = 50
$ws1.Cells.Item(4,8).Value = 'https://show.bilibili.com/platform/detail.html?id=84925'
$ws1.Cells.Item(4,9).Value = '//i2.hdslb.com/bfs/openplatform/202404/YjFyyYq51713508727131.jpeg'
$ws1.Cells.Item(5,2).NumberFormat = "@"
$ws1.Cells.Item(5,2).Value = '2024-07-06'
$ws1.Cells.Item(5,2).Style = "Normal"
$ws1.Cells.Item(5,3).Value = '南宁·首届童话梦境Lolita茶会'
$ws1.Cells.Item(5,4).Value = '明秀东路157号 利泰国际大酒店'
$ws1.Cells.Item(5,5).Value = '2024.07.06 13:00-07.06 17:00'
$ws1.Cells.Item(5,6).Value = 144
$ws1.Cells.Item(5,7).Value = 88
$ws1.Cells.Item(5,8).Value = 'https://show.bilibili.com/platform/detail.html?id=85776'
$ws1.Cells.Item(5,9).Value = '//i2.hdslb.com/bfs/openplatform/202405/Xl4NBnky1715847180514.jpeg'
$ws1.Cells.Item(6,2).NumberFormat = "@"
$ws1.Cells.Item(6,2).Value = '2024-07-13'
$ws1.Cells.Item(6,2).Style = "Normal"
$ws1.Cells.Item(6,3).Value = '南宁·0713国乙ONLY'
$ws1.Cells.Item(6,4).Value = '亭洪路45号 水明漾宴会中心'
$ws1.Cells.Item(6,5).Value = '2024.07.13 09:30-07.13 21:00'
$ws1.Cells.Item(6,6).Value = 229
$ws1.Cells.Item(6,7).Value = 68
$ws1.Cells.Item(6,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86378'
$ws1.Cells.Item(6,9).Value = '//i1.hdslb.com/bfs/openplatform/202405/ZDBCv2of1716659486569.jpeg'
$ws1.Cells.Item(7,2).NumberFormat = "@"
$ws1.Cells.Item(7,2).Value = '2024-07-14'
$ws1.Cells.Item(7,2).Style = "Normal"
$ws1.Cells.Item(7,3).Value = '广西·首届明日方舟only展 - 花庭圣梦'
$ws1.Cells.Item(7,4).Value = '明秀东路157号 利泰国际大酒店'
$ws1.Cells.Item(7,5).Value = '2024.07.14 09:00-07.14 18:00'
$ws1.Cells.Item(7,6).Value = 190
$ws1.Cells.Item(7,7).Value = 69
$ws1.Cells.Item(7,8).Value = 'https://show.bilibili.com/platform/detail.html?id=85852'
$ws1.Cells.Item(7,9).Value = '//i2.hdslb.com/bfs/openplatform/202405/xsMTmueN1715920435584.jpeg'
$ws1.Cells.Item(8,2).NumberFormat = "@"
$ws1.Cells.Item(8,2).Value = '2024-07-20'
$ws1.Cells.Item(8,2).Style = "Normal"
$ws1.Cells.Item(8,3).Value = '南宁·AB动漫游戏嘉年华'
$ws1.Cells.Item(8,4).Value = '三塘南路与长虹东路交叉路口往北约50米 广西农业会展中心'
$ws1.Cells.Item(8,5).Value = '2024.07.20 09:30-07.21 17:00'
$ws1.Cells.Item(8,6).Value = 1839
$ws1.Cells.Item(8,7).Value = 60
$ws1.Cells.Item(8,8).Value = 'https://show.bilibili.com/platform/detail.html?id=84862'
$ws1.Cells.Item(8,9).Value = '//i1.hdslb.com/bfs/openplatform/202404/eglavDeZ1714036487217.jpeg'
$ws1.Cells.Item(9,2).NumberFormat = "@"
$ws1.Cells.Item(9,2).Value = '2024-07-20'
$ws1.Cells.Item(9,2).Style = "Normal"
$ws1.Cells.Item(9,3).Value = '横州·第二届海棠动漫游戏嘉年华'
$ws1.Cells.Item(9,4).Value = '茉莉花大道 横州国际大酒店'
$ws1.Cells.Item(9,5).Value = '2024.07.20 09:30-07.20 17:00'
$ws1.Cells.Item(9,6).Value = 330
$ws1.Cells.Item(9,7).Value = 30
$ws1.Cells.Item(9,8).Value = 'https://show.bilibili.com/platform/detail.html?id=84799'
$ws1.Cells.Item(9,9).Value = '//i2.hdslb.com/bfs/openplatform/202404/r50S2ttT1713869164413.jpeg'
$ws1.Cells.Item(10,2).NumberFormat = "@"
$ws1.Cells.Item(10,2).Value = '2024-07-27'
$ws1.Cells.Item(10,2).Style = "Normal"
$ws1.Cells.Item(10,3).Value = '南宁·第十九届（2024）良牙动漫夏季盛典（良牙夏典）'
$ws1.Cells.Item(10,4).Value = '民族大道106号 南宁国际会展中心'
$ws1.Cells.Item(10,5).Value = '2024.07.27 09:30-07.28 17:30'
$ws1.Cells.Item(10,6).Value = 4274
$ws1.Cells.Item(10,7).Value = 55
$ws1.Cells.Item(10,8).Value = 'https://show.bilibili.com/platform/detail.html?id=85264'
$ws1.Cells.Item(10,9).Value = '//i0.hdslb.com/bfs/openplatform/202405/dZVcS7eE1715155418142.jpeg'
$ws1.Cells.Item(11,2).NumberFormat = "@"
$ws1.Cells.Item(11,2).Value = '2024-08-03'
$ws1.Cells.Item(11,2).Style = "Normal"
$ws1.Cells.Item(11,3).Value = '南宁·火影忍者only'
$ws1.Cells.Item(11,4).Value = '厢竹大道65号 桔子酒店'
$ws1.Cells.Item(11,5).Value = '2024.08.03 10:00-08.03 17:00'
$ws1.Cells.Item(11,6).Value = 54
$ws1.Cells.Item(11,7).Value = 68
$ws1.Cells.Item(11,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86994'
$ws1.Cells.Item(11,9).Value = '//i0.hdslb.com/bfs/openplatform/202406/h1tXE9t11717523356034.jpeg'
$ws1.Cells.Item(12,2).NumberFormat = "@"
$ws1.Cells.Item(12,2).Value = '2024-08-03'
$ws1.Cells.Item(12,2).Style = "Normal"
$ws1.Cells.Item(12,3).Value = '南宁·蔚蓝档案only'
$ws1.Cells.Item(12,4).Value = '亭洪路45号 百益上河城'
$ws1.Cells.Item(12,5).Value = '2024.08.03 09:00-08.03 17:00'
$ws1.Cells.Item(12,6).Value = 295
$ws1.Cells.Item(12,7).Value = 68
$ws1.Cells.Item(12,8).Value = 'https://show.bilibili.com/platform/detail.html?id=85370'
$ws1.Cells.Item(12,9).Value = '//i1.hdslb.com/bfs/openplatform/202405/sBxi2Mx61715247424836.jpeg'

# ---------------------------------------------------------------------------
# Sheet "全部类型" (4th sheet): 18 data+header rows -> 16 data+header rows.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

# Drop the oldest record (old row 2) -- everything below shifts up one row.
$ws4.Rows.Item(2).Delete()
# Drop what is now the trailing record (old row 18, now row 17).
$ws4.Rows.Item(17).Delete()

# Renumber the "0"-based index column (A) so it stays sequential (1..15).
for ($r = 2; $r -le 16; $r++) {
    $ws4.Cells.Item($r, 1).Value = $r - 1
}

# Refresh every remaining row's scraped fields (B:I) to the latest values.
$ws4.Cells.Item(2,2).NumberFormat = "@"
$ws4.Cells.Item(2,2).Value = '2024-06-15'
$ws4.Cells.Item(2,2).Style = "Normal"
$ws4.Cells.Item(2,3).Value = '南宁·《菊次郎的夏天》久石让作品视听音乐会'
$ws4.Cells.Item(2,4).Value = '民族大道49-2号 广西音乐厅（广西民族艺术宫）'
$ws4.Cells.Item(2,5).Value = '2024.06.15 20:00-06.15 21:30'
$ws4.Cells.Item(2,6).Value = 1
$ws4.Cells.Item(2,7).Value = 108
$ws4.Cells.Item(2,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86653'
$ws4.Cells.Item(2,9).Value = '//i0.hdslb.com/bfs/openplatform/202405/v2g0hMrK1717123700770.png'
$ws4.Cells.Item(3,2).NumberFormat = "@"
$ws4.Cells.Item(3,2).Value = '2024-06-15'
$ws4.Cells.Item(3,2).Style = "Normal"
$ws4.Cells.Item(3,3).Value = '南宁·星STAR国潮嘉年华'
$ws4.Cells.Item(3,4).Value = '亭洪路45号 百益上河城'
$ws4.Cells.Item(3,5).Value = '2024.06.15 09:00-06.16 17:00'
$ws4.Cells.Item(3,6).Value = 73
$ws4.Cells.Item(3,7).Value = 50
$ws4.Cells.Item(3,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86198'
$ws4.Cells.Item(3,9).Value = '//i0.hdslb.com/bfs/openplatform/202405/orwMgait1716448294056.jpeg'
$ws4.Cells.Item(4,2).NumberFormat = "@"
$ws4.Cells.Item(4,2).Value = '2024-06-22'
$ws4.Cells.Item(4,2).Style = "Normal"
$ws4.Cells.Item(4,3).Value = '南宁·排球少年ONLY（取消）'
$ws4.Cells.Item(4,4).Value = '亭洪路45号 水明漾宴会中心'
$ws4.Cells.Item(4,5).Value = '2024.06.22 09:45-06.22 17:00'
$ws4.Cells.Item(4,6).Value = 63
$ws4.Cells.Item(4,7).Value = '不可售'
$ws4.Cells.Item(4,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86465'
$ws4.Cells.Item(4,9).Value = '//i0.hdslb.com/bfs/openplatform/202405/GaaD97dL1716883956953.jpeg'
$ws4.Cells.Item(5,2).NumberFormat = "@"
$ws4.Cells.Item(5,2).Value = '2024-06-22'
$ws4.Cells.Item(5,2).Style = "Normal"
$ws4.Cells.Item(5,3).Value = '南宁·浪漫古典·百年经典世界名曲音乐会'
$ws4.Cells.Item(5,4).Value = '广西壮族自治区南宁市良庆区龙堤路25号  广西文化艺术中心-音乐厅'
$ws4.Cells.Item(5,5).Value = '2024.06.22 20:00-06.22 21:30'
$ws4.Cells.Item(5,6).Value = 48
$ws4.Cells.Item(5,7).Value = 50
$ws4.Cells.Item(5,8).Value = 'https://show.bilibili.com/platform/detail.html?id=83959'
$ws4.Cells.Item(5,9).Value = '//i1.hdslb.com/bfs/openplatform/202404/H0f8U7no1712041461015.jpeg'
$ws4.Cells.Item(6,2).NumberFormat = "@"
$ws4.Cells.Item(6,2).Value = '2024-07-06'
$ws4.Cells.Item(6,2).Style = "Normal"
$ws4.Cells.Item(6,3).Value = '南宁·小蜜蜂动漫嘉年华2.0'
$ws4.Cells.Item(6,4).Value = '亭洪路45号 百益上河城'
$ws4.Cells.Item(6,5).Value = '2024.07.06 10:00-07.06 17:00'
$ws4.Cells.Item(6,6).Value = 245
$ws4.Cells.Item(6,7).Value = 50
$ws4.Cells.Item(6,8).Value = 'https://show.bilibili.com/platform/detail.html?id=84925'
$ws4.Cells.Item(6,9).Value = '//i2.hdslb.com/bfs/openplatform/202404/YjFyyYq51713508727131.jpeg'
$ws4.Cells.Item(7,2).NumberFormat = "@"
$ws4.Cells.Item(7,2).Value = '2024-07-06'
$ws4.Cells.Item(7,2).Style = "Normal"
$ws4.Cells.Item(7,3).Value = '南宁·首届童话梦境Lolita茶会'
$ws4.Cells.Item(7,4).Value = '明秀东路157号 利泰国际大酒店'
$ws4.Cells.Item(7,5).Value = '2024.07.06 13:00-07.06 17:00'
$ws4.Cells.Item(7,6).Value = 144
$ws4.Cells.Item(7,7).Value = 88
$ws4.Cells.Item(7,8).Value = 'https://show.bilibili.com/platform/detail.html?id=85776'
$ws4.Cells.Item(7,9).Value = '//i2.hdslb.com/bfs/openplatform/202405/Xl4NBnky1715847180514.jpeg'
$ws4.Cells.Item(8,2).NumberFormat = "@"
$ws4.Cells.Item(8,2).Value = '2024-07-13'
$ws4.Cells.Item(8,2).Style = "Normal"
$ws4.Cells.Item(8,3).Value = '南宁·0713国乙ONLY'
$ws4.Cells.Item(8,4).Value = '亭洪路45号 水明漾宴会中心'
$ws4.Cells.Item(8,5).Value = '2024.07.13 09:30-07.13 21:00'
$ws4.Cells.Item(8,6).Value = 229
$ws4.Cells.Item(8,7).Value = 68
$ws4.Cells.Item(8,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86378'
$ws4.Cells.Item(8,9).Value = '//i1.hdslb.com/bfs/openplatform/202405/ZDBCv2of1716659486569.jpeg'
$ws4.Cells.Item(9,2).NumberFormat = "@"
$ws4.Cells.Item(9,2).Value = '2024-07-14'
$ws4.Cells.Item(9,2).Style = "Normal"
$ws4.Cells.Item(9,3).Value = '广西·首届明日方舟only展 - 花庭圣梦'
$ws4.Cells.Item(9,4).Value = '明秀东路157号 利泰国际大酒店'
$ws4.Cells.Item(9,5).Value = '2024.07.14 09:00-07.14 18:00'
$ws4.Cells.Item(9,6).Value = 190
$ws4.Cells.Item(9,7).Value = 69
$ws4.Cells.Item(9,8).Value = 'https://show.bilibili.com/platform/detail.html?id=85852'
$ws4.Cells.Item(9,9).Value = '//i2.hdslb.com/bfs/openplatform/202405/xsMTmueN1715920435584.jpeg'
$ws4.Cells.Item(10,2).NumberFormat = "@"
$ws4.Cells.Item(10,2).Value = '2024-07-18'
$ws4.Cells.Item(10,2).Style = "Normal"
$ws4.Cells.Item(10,3).Value = '南宁·限时6折|俄罗斯圣彼得堡古典芭蕾舞剧院《胡桃夹子》'
$ws4.Cells.Item(10,4).Value = '龙堤路25号 广西文化艺术中心'
$ws4.Cells.Item(10,5).Value = '2024.07.18 20:00-07.18 21:30'
$ws4.Cells.Item(10,6).Value = 5
$ws4.Cells.Item(10,7).Value = 108
$ws4.Cells.Item(10,8).Value = 'https://show.bilibili.com/platform/detail.html?id=85816'
$ws4.Cells.Item(10,9).Value = '//i0.hdslb.com/bfs/openplatform/202405/SN0ZyGVj1715675672714.jpeg'
$ws4.Cells.Item(11,2).NumberFormat = "@"
$ws4.Cells.Item(11,2).Value = '2024-07-19'
$ws4.Cells.Item(11,2).Style = "Normal"
$ws4.Cells.Item(11,3).Value = '南宁·限时6折|俄罗斯圣彼得堡古典芭蕾舞剧院《天鹅湖》 '
$ws4.Cells.Item(11,4).Value = '龙堤路25号 广西文化艺术中心'
$ws4.Cells.Item(11,5).Value = '2024.07.19 20:00-07.19 22:00'
$ws4.Cells.Item(11,6).Value = 10
$ws4.Cells.Item(11,7).Value = 108
$ws4.Cells.Item(11,8).Value = 'https://show.bilibili.com/platform/detail.html?id=85831'
$ws4.Cells.Item(11,9).Value = '//i1.hdslb.com/bfs/openplatform/202405/ZyyeeOUo1715677877362.jpeg'
$ws4.Cells.Item(12,2).NumberFormat = "@"
$ws4.Cells.Item(12,2).Value = '2024-07-20'
$ws4.Cells.Item(12,2).Style = "Normal"
$ws4.Cells.Item(12,3).Value = '南宁·AB动漫游戏嘉年华'
$ws4.Cells.Item(12,4).Value = '三塘南路与长虹东路交叉路口往北约50米 广西农业会展中心'
$ws4.Cells.Item(12,5).Value = '2024.07.20 09:30-07.21 17:00'
$ws4.Cells.Item(12,6).Value = 1839
$ws4.Cells.Item(12,7).Value = 60
$ws4.Cells.Item(12,8).Value = 'https://show.bilibili.com/platform/detail.html?id=84862'
$ws4.Cells.Item(12,9).Value = '//i1.hdslb.com/bfs/openplatform/202404/eglavDeZ1714036487217.jpeg'
$ws4.Cells.Item(13,2).NumberFormat = "@"
$ws4.Cells.Item(13,2).Value = '2024-07-20'
$ws4.Cells.Item(13,2).Style = "Normal"
$ws4.Cells.Item(13,3).Value = '横州·第二届海棠动漫游戏嘉年华'
$ws4.Cells.Item(13,4).Value = '茉莉花大道 横州国际大酒店'
$ws4.Cells.Item(13,5).Value = '2024.07.20 09:30-07.20 17:00'
$ws4.Cells.Item(13,6).Value = 330
$ws4.Cells.Item(13,7).Value = 30
$ws4.Cells.Item(13,8).Value = 'https://show.bilibili.com/platform/detail.html?id=84799'
$ws4.Cells.Item(13,9).Value = '//i2.hdslb.com/bfs/openplatform/202404/r50S2ttT1713869164413.jpeg'
$ws4.Cells.Item(14,2).NumberFormat = "@"
$ws4.Cells.Item(14,2).Value = '2024-07-27'
$ws4.Cells.Item(14,2).Style = "Normal"
$ws4.Cells.Item(14,3).Value = '南宁·第十九届（2024）良牙动漫夏季盛典（良牙夏典）'
$ws4.Cells.Item(14,4).Value = '民族大道106号 南宁国际会展中心'
$ws4.Cells.Item(14,5).Value = '2024.07.27 09:30-07.28 17:30'
$ws4.Cells.Item(14,6).Value = 4274
$ws4.Cells.Item(14,7).Value = 55
$ws4.Cells.Item(14,8).Value = 'https://show.bilibili.com/platform/detail.html?id=85264'
$ws4.Cells.Item(14,9).Value = '//i0.hdslb.com/bfs/openplatform/202405/dZVcS7eE1715155418142.jpeg'
$ws4.Cells.Item(15,2).NumberFormat = "@"
$ws4.Cells.Item(15,2).Value = '2024-08-03'
$ws4.Cells.Item(15,2).Style = "Normal"
$ws4.Cells.Item(15,3).Value = '南宁·火影忍者only'
$ws4.Cells.Item(15,4).Value = '厢竹大道65号 桔子酒店'
$ws4.Cells.Item(15,5).Value = '2024.08.03 10:00-08.03 17:00'
$ws4.Cells.Item(15,6).Value = 54
$ws4.Cells.Item(15,7).Value = 68
$ws4.Cells.Item(15,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86994'
$ws4.Cells.Item(15,9).Value = '//i0.hdslb.com/bfs/openplatform/202406/h1tXE9t11717523356034.jpeg'
$ws4.Cells.Item(16,2).NumberFormat = "@"
$ws4.Cells.Item(16,2).Value = '2024-08-03'
$ws4.Cells.Item(16,2).Style = "Normal"
$ws4.Cells.Item(16,3).Value = '南宁·蔚蓝档案only'
$ws4.Cells.Item(16,4).Value = '亭洪路45号 百益上河城'
$ws4.Cells.Item(16,5).Value = '2024.08.03 09:00-08.03 17:00'
$ws4.Cells.Item(16,6).Value = 295
$ws4.Cells.Item(16,7).Value = 68
$ws4.Cells.Item(16,8).Value = 'https://show.bilibili.com/platform/detail.html?id=85370'
$ws4.Cells.Item(16,9).Value = '//i1.hdslb.com/bfs/openplatform/202405/sBxi2Mx61715247424836.jpeg'
